$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.803092333333333
$ws.Range("H2").Value = 5.409276999999999
$ws.Range("I2").Value = 0.1744886524959502
$ws.Range("J2").Value = 0.1744886524959502
$ws.Range("M2").Value = 8.252454666666667
$ws.Range("N2").Value = 24.757364
$ws.Range("O2").Value = 0.05349680956196952
$ws.Range("P2").Value = 0.05349680956196953
$ws.Range("Q2").Value = 14.87993774064755
$ws.Range("R2").Value = 133.919439665828
$ws.Range("S2").Value = 0.009334586213300527
$ws.Range("T2").Value = 0.009334586213300527
$ws.Range("G3").Value = 1.803092333333333
$ws.Range("H3").Value = 5.409276999999999
$ws.Range("I3").Value = 0.1744886524959502
$ws.Range("J3").Value = 0.1744886524959502
$ws.Range("O3").Value = 0.5638948237978928
$ws.Range("P3").Value = 0.5638948237978929
$ws.Range("Q3").Value = 156.8452387925385
$ws.Range("R3").Value = 1411.607149132847
$ws.Range("S3").Value = 0.09839324795393561
$ws.Range("T3").Value = 0.09839324795393561
$ws.Range("G4").Value = 1.803092333333333
$ws.Range("H4").Value = 5.409276999999999
$ws.Range("I4").Value = 0.1744886524959502
$ws.Range("J4").Value = 0.1744886524959502
$ws.Range("M4").Value = 57.81408433333333
$ws.Range("N4").Value = 173.442253
$ws.Range("O4").Value = 0.3747817085348802
$ws.Range("P4").Value = 0.3747817085348802
$ws.Range("Q4").Value = 104.2441322201201
$ws.Range("R4").Value = 938.1971899810809
$ws.Range("S4").Value = 0.06539515530238121
$ws.Range("T4").Value = 0.06539515530238121
$ws.Range("G5").Value = 1.803092333333333
$ws.Range("H5").Value = 5.409276999999999
$ws.Range("I5").Value = 0.1744886524959502
$ws.Range("J5").Value = 0.1744886524959502
$ws.Range("M5").Value = 1.207345666666667
$ws.Range("N5").Value = 3.622037
$ws.Range("O5").Value = 0.007826658105257385
$ws.Range("P5").Value = 0.007826658105257386
$ws.Range("Q5").Value = 2.176955715249889
$ws.Range("R5").Value = 19.592601437249
$ws.Range("S5").Value = 0.001365663026332868
$ws.Range("T5").Value = 0.001365663026332868
$ws.Range("I6").Value = 0.4384883998568034
$ws.Range("J6").Value = 0.4384883998568034
$ws.Range("M6").Value = 8.252454666666667
$ws.Range("N6").Value = 24.757364
$ws.Range("O6").Value = 0.05349680956196952
$ws.Range("P6").Value = 0.05349680956196953
$ws.Range("Q6").Value = 37.39314847432178
$ws.Range("R6").Value = 336.538336268896
$ws.Range("S6").Value = 0.02345773042227215
$ws.Range("T6").Value = 0.02345773042227215
$ws.Range("I7").Value = 0.4384883998568034
$ws.Range("J7").Value = 0.4384883998568034
$ws.Range("O7").Value = 0.5638948237978928
$ws.Range("P7").Value = 0.5638948237978929
$ws.Range("R7").Value = 3547.355952353704
$ws.Range("S7").Value = 0.2472613389746721
$ws.Range("T7").Value = 0.2472613389746721
$ws.Range("I8").Value = 0.4384883998568034
$ws.Range("J8").Value = 0.4384883998568034
$ws.Range("M8").Value = 57.81408433333333
$ws.Range("N8").Value = 173.442253
$ws.Range("O8").Value = 0.3747817085348802
$ws.Range("P8").Value = 0.3747817085348802
$ws.Range("Q8").Value = 261.9645580260436
$ws.Range("R8").Value = 2357.681022234392
$ws.Range("S8").Value = 0.1643374316710585
$ws.Range("T8").Value = 0.1643374316710585
$ws.Range("I9").Value = 0.4384883998568034
$ws.Range("J9").Value = 0.4384883998568034
$ws.Range("M9").Value = 1.207345666666667
$ws.Range("N9").Value = 3.622037
$ws.Range("O9").Value = 0.007826658105257385
$ws.Range("P9").Value = 0.007826658105257386
$ws.Range("Q9").Value = 5.470669951796444
$ws.Range("R9").Value = 49.236029566168
$ws.Range("S9").Value = 0.003431898788800592
$ws.Range("T9").Value = 0.003431898788800592
$ws.Range("G10").Value = 3.895605666666667
$ws.Range("H10").Value = 11.686817
$ws.Range("I10").Value = 0.3769851220961256
$ws.Range("J10").Value = 0.3769851220961256
$ws.Range("M10").Value = 8.252454666666667
$ws.Range("N10").Value = 24.757364
$ws.Range("O10").Value = 0.05349680956196952
$ws.Range("P10").Value = 0.05349680956196953
$ws.Range("Q10").Value = 32.14830916337645
$ws.Range("R10").Value = 289.3347824703881
$ws.Range("S10").Value = 0.02016750128447226
$ws.Range("T10").Value = 0.02016750128447226
$ws.Range("G11").Value = 3.895605666666667
$ws.Range("H11").Value = 11.686817
$ws.Range("I11").Value = 0.3769851220961256
$ws.Range("J11").Value = 0.3769851220961256
$ws.Range("O11").Value = 0.5638948237978928
$ws.Range("P11").Value = 0.5638948237978929
$ws.Range("Q11").Value = 338.8662852890875
$ws.Range("R11").Value = 3049.796567601788
$ws.Range("S11").Value = 0.2125799589988218
$ws.Range("T11").Value = 0.2125799589988219
$ws.Range("G12").Value = 3.895605666666667
$ws.Range("H12").Value = 11.686817
$ws.Range("I12").Value = 0.3769851220961256
$ws.Range("J12").Value = 0.3769851220961256
$ws.Range("M12").Value = 57.81408433333333
$ws.Range("N12").Value = 173.442253
$ws.Range("O12").Value = 0.3747817085348802
$ws.Range("P12").Value = 0.3747817085348802
$ws.Range("Q12").Value = 225.2208745420779
$ws.Range("R12").Value = 2026.987870878701
$ws.Range("S12").Value = 0.1412871281514164
$ws.Range("T12").Value = 0.1412871281514164
$ws.Range("G13").Value = 3.895605666666667
$ws.Range("H13").Value = 11.686817
$ws.Range("I13").Value = 0.3769851220961256
$ws.Range("J13").Value = 0.3769851220961256
$ws.Range("M13").Value = 1.207345666666667
$ws.Range("N13").Value = 3.622037
$ws.Range("O13").Value = 0.007826658105257385
$ws.Range("P13").Value = 0.007826658105257386
$ws.Range("Q13").Value = 4.703342620692111
$ws.Range("R13").Value = 42.33008358622901
$ws.Range("S13").Value = 0.002950533661415086
$ws.Range("T13").Value = 0.002950533661415087
$ws.Range("G14").Value = 0.1037266666666667
$ws.Range("H14").Value = 0.31118
$ws.Range("I14").Value = 0.01003782555112075
$ws.Range("J14").Value = 0.01003782555112075
$ws.Range("M14").Value = 8.252454666666667
$ws.Range("N14").Value = 24.757364
$ws.Range("O14").Value = 0.05349680956196952
$ws.Range("P14").Value = 0.05349680956196953
$ws.Range("Q14").Value = 0.8559996143911112
$ws.Range("R14").Value = 7.703996529520001
$ws.Range("S14").Value = 0.0005369916419245786
$ws.Range("T14").Value = 0.0005369916419245786
$ws.Range("G15").Value = 0.1037266666666667
$ws.Range("H15").Value = 0.31118
$ws.Range("I15").Value = 0.01003782555112075
$ws.Range("J15").Value = 0.01003782555112075
$ws.Range("O15").Value = 0.5638948237978928
$ws.Range("P15").Value = 0.5638948237978929
$ws.Range("Q15").Value = 9.022851188331112
$ws.Range("R15").Value = 81.20566069498001
$ws.Range("S15").Value = 0.005660277870463222
$ws.Range("T15").Value = 0.005660277870463222
$ws.Range("G16").Value = 0.1037266666666667
$ws.Range("H16").Value = 0.31118
$ws.Range("I16").Value = 0.01003782555112075
$ws.Range("J16").Value = 0.01003782555112075
$ws.Range("M16").Value = 57.81408433333333
$ws.Range("N16").Value = 173.442253
$ws.Range("O16").Value = 0.3747817085348802
$ws.Range("P16").Value = 0.3747817085348802
$ws.Range("Q16").Value = 5.996862254282223
$ws.Range("R16").Value = 53.97176028854
$ws.Range("S16").Value = 0.003761993410024111
$ws.Range("T16").Value = 0.003761993410024111
$ws.Range("G17").Value = 0.1037266666666667
$ws.Range("H17").Value = 0.31118
$ws.Range("I17").Value = 0.01003782555112075
$ws.Range("J17").Value = 0.01003782555112075
$ws.Range("M17").Value = 1.207345666666667
$ws.Range("N17").Value = 3.622037
$ws.Range("O17").Value = 0.007826658105257385
$ws.Range("P17").Value = 0.007826658105257386
$ws.Range("Q17").Value = 0.1252339415177778
$ws.Range("R17").Value = 1.12710547366
$ws.Range("S17").Value = 0.00007856262870883891
$ws.Range("T17").Value = 0.00007856262870883891
